$d = $word.ActiveDocument

# --- Locate the block of paragraphs to replace -----------------------------
# Start: the paragraph that currently reads "Add features for phone (Descriptions)"
# End:   the paragraph that holds the "_GoBack" bookmark (last paragraph of the
#        block, right before the trailing empty paragraph we also need to add).
$startFind = $d.Content
$startFind.Find.Execute("Add features for phone (Descriptions)") | Out-Null
$startPos = $startFind.Paragraphs(1).Range.Start

$bm = $d.Bookmarks("_GoBack")
$endPos = $bm.Range.Paragraphs(1).Range.End

$targetRange = $d.Range($startPos, $endPos)

# --- Build the replacement OOXML --------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = ''
$newXml += '<w:p ' + $wNs + '>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>'
$newXml +=   '<w:r><w:lastRenderedPageBreak/><w:t>2/4/18</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$newXml +=   '<w:r><w:t>Implemented majority of hardware (one sensor not working)</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$newXml +=   '<w:r><w:t>Added prepared statements for SQL</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$newXml +=   '<w:r><w:t>Fixed quiz so it asks different questions every time</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>'
$newXml +=   '<w:r><w:t>09-4-18</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>'
$newXml +=   '<w:r><w:t>Began to finish design of quiz.</w:t></w:r>'
$newXml += '</w:p>'

$newXml += '<w:p>'
$newXml +=   '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>'
$newXml +=   '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$newXml += '</w:p>'

$newXml += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'

$targetRange.InsertXML($newXml)
